{"js": "// The diff inserts two new paragraphs (\"Titel\" style) at the very top of\n// the document body, immediately before the existing first paragraph:\n//   1) an empty paragraph\n//   2) a paragraph with the text \"Eine dritte \u00e4nderung local am word\"\nconst body = context.document.body;\nconst firstPara = body.paragraphs.getFirst();\n\n// Insert the new titled paragraph right before the current first paragraph.\nconst textPara = firstPara.insertParagraph(\n  \"Eine dritte \u00e4nderung local am word\",\n  Word.InsertLocation.before\n);\ntextPara.styleBuiltIn = Word.Style.title;\n\n// Insert an empty paragraph right before that new text paragraph, so the\n// final order is: [empty Titel] -> [Titel \"Eine dritte \u00e4nderung local am\n// word\"] -> [original first paragraph] -> ...rest of the document.\nconst emptyPara = textPara.insertParagraph(\"\", Word.InsertLocation.before);\nemptyPara.styleBuiltIn = Word.Style.title;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Work against the very first paragraph of the document (currently the\n# \"\u00c4nderung zum 2.\" title paragraph).\n$firstPara = $d.Paragraphs(1)\n$r = $firstPara.Range\n\n# Insert two new paragraphs immediately before it (this shifts the\n# original first paragraph down by two, and both new paragraphs inherit\n# the \"Titel\" style from the paragraph they were split from).\n$r.InsertParagraphBefore()\n$r.InsertParagraphBefore()\n\n# Paragraph 1 stays empty (just the \"Titel\"-styled paragraph mark).\n# Paragraph 2 gets the new title text.\n$d.Paragraphs(2).Range.Text = \"Eine dritte \u00e4nderung local am word\"\n\n# Make sure both inserted paragraphs explicitly carry the \"Titel\" style.\n$d.Paragraphs(1).Style = \"Titel\"\n$d.Paragraphs(2).Style = \"Titel\"\n"}
